# Rename "Guess" -> "Gaussian", reorder its parameter-guess rows, and add a
# new "Lorentzian" sheet (with its own parameter guesses) right after it.

$wb = $excel.ActiveWorkbook

# --- Rename Guess -> Gaussian -----------------------------------------
$gaussian = $wb.Worksheets.Item("Guess")
$gaussian.Name = "Gaussian"

# Reorder rows 3-5 on the Gaussian sheet to: Abase, sigma, Agaussian
$gaussian.Cells.Item(3, 1).Value = "Abase"
$gaussian.Cells.Item(3, 2).Value = 0.9
$gaussian.Cells.Item(3, 3).Value = 0.91
$gaussian.Cells.Item(3, 4).Value = 0.93
$gaussian.Cells.Item(3, 5).Value = 0.9
$gaussian.Cells.Item(3, 6).Value = 0.82
$gaussian.Cells.Item(3, 7).Value = 0.45

$gaussian.Cells.Item(4, 1).Value = "sigma"
$gaussian.Cells.Item(4, 2).Value = 1.3
$gaussian.Cells.Item(4, 3).Value = 1.3
$gaussian.Cells.Item(4, 4).Value = 1.3
$gaussian.Cells.Item(4, 5).Value = 1.3
$gaussian.Cells.Item(4, 6).Value = 2
$gaussian.Cells.Item(4, 7).Value = 1.59

$gaussian.Cells.Item(5, 1).Value = "Agaussian"
$gaussian.Cells.Item(5, 2).Value = 4000
$gaussian.Cells.Item(5, 3).Value = 3250
$gaussian.Cells.Item(5, 4).Value = 2250
$gaussian.Cells.Item(5, 5).Value = 2250
$gaussian.Cells.Item(5, 6).Value = 2250
$gaussian.Cells.Item(5, 7).Value = 9250

[void]$gaussian.Activate()
[void]$gaussian.Range("A4:H5").Select()

# --- Add the Lorentzian sheet right after Gaussian --------------------
$lorentzian = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $gaussian)
$lorentzian.Name = "Lorentzian"

$lorentzian.Cells.Item(1, 1).Value = "Variables"
$lorentzian.Cells.Item(1, 2).Value = "ann1"
$lorentzian.Cells.Item(1, 3).Value = "pt2"
$lorentzian.Cells.Item(1, 4).Value = "pt2b"
$lorentzian.Cells.Item(1, 5).Value = "pt2c"
$lorentzian.Cells.Item(1, 6).Value = "pt2d"
$lorentzian.Cells.Item(1, 7).Value = "pt2e"

$lorentzian.Cells.Item(2, 1).Value = "x0"
$lorentzian.Cells.Item(2, 2).Value = 0
$lorentzian.Cells.Item(2, 3).Value = 0
$lorentzian.Cells.Item(2, 4).Value = 0
$lorentzian.Cells.Item(2, 5).Value = 0
$lorentzian.Cells.Item(2, 6).Value = 0
$lorentzian.Cells.Item(2, 7).Value = 0

$lorentzian.Cells.Item(3, 1).Value = "Abase"
$lorentzian.Cells.Item(3, 2).Value = 0.9
$lorentzian.Cells.Item(3, 3).Value = 0.9
$lorentzian.Cells.Item(3, 4).Value = 0.9
$lorentzian.Cells.Item(3, 5).Value = 0.9
$lorentzian.Cells.Item(3, 6).Value = 0.85
$lorentzian.Cells.Item(3, 7).Value = 0.85

$lorentzian.Cells.Item(4, 1).Value = "gamma"
$lorentzian.Cells.Item(4, 2).Value = 0.8
$lorentzian.Cells.Item(4, 3).Value = 0.8
$lorentzian.Cells.Item(4, 4).Value = 0.8
$lorentzian.Cells.Item(4, 5).Value = 0.8
$lorentzian.Cells.Item(4, 6).Value = 0.8
$lorentzian.Cells.Item(4, 7).Value = 0.8

$lorentzian.Cells.Item(5, 1).Value = "Alorentzian"
$lorentzian.Cells.Item(5, 2).Value = 3500
$lorentzian.Cells.Item(5, 3).Value = 3500
$lorentzian.Cells.Item(5, 4).Value = 3500
$lorentzian.Cells.Item(5, 5).Value = 3500
$lorentzian.Cells.Item(5, 6).Value = 2000
$lorentzian.Cells.Item(5, 7).Value = 2000

[void]$lorentzian.Activate()
[void]$lorentzian.Range("A4:G5").Select()
